$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 74, shifting existing rows 74:179 down to 75:180
# (this also extends the used range / dimension to A1:R180 automatically,
# and carries the existing row formatting, e.g. the date style on column D).
$ws.Rows("74:74").Insert()

# Fill in the new row 74 with the new data record.
# Columns A,B,C,E,F,G,H,I,O,R keep the same constant values used throughout
# this sheet (single market / category subset), matching the diff context.
$ws.Range("A74").Value = 4
$ws.Range("B74").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C74").Value = "Los Lagos"
$ws.Range("D74").Value = 44579
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = 100112032
$ws.Range("G74").Value = "Zapallo italiano"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 200
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 16000
$ws.Range("M74").Value = 15600
$ws.Range("N74").Value = "$/caja 50 unidades"
$ws.Range("O74").Value = "Región de O'Higgins"
$ws.Range("P74").Value = 312
$ws.Range("Q74").Value = 50
$ws.Range("R74").Value = "Hortaliza"
